$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing D:K data to F:M
$ws.Columns("D:E").Insert()

# Copy number formats/styles from the (now-shifted) F:G columns back onto the
# new D:E columns so the new quarters look like the rest of the table.
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarters of data
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 181100
$ws.Range("E8").Value2 = 187700
$ws.Range("D9").Value2 = 49400
$ws.Range("E9").Value2 = 54200
$ws.Range("D10").Value2 = 131700
$ws.Range("E10").Value2 = 133500
$ws.Range("D17").Value2 = 200500
$ws.Range("E17").Value2 = 203500
$ws.Range("D18").Value2 = -19400
$ws.Range("E18").Value2 = -15800
$ws.Range("D21").Value2 = -17100
$ws.Range("E21").Value2 = -12600
$ws.Range("D23").Value2 = -19400
$ws.Range("E23").Value2 = -15900
$ws.Range("D24").Value2 = -8400
$ws.Range("E24").Value2 = -4200
$ws.Range("D26").Value2 = -11000
$ws.Range("E26").Value2 = -11700
$ws.Range("D27").Value2 = -11100
$ws.Range("E27").Value2 = -11700
$ws.Range("D33").Value2 = -11100
$ws.Range("E33").Value2 = -11700
$ws.Range("D35").Value2 = -11100
$ws.Range("E35").Value2 = -11700
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 112700
$ws.Range("E41").Value2 = 241000
$ws.Range("D43").Value2 = 95800
$ws.Range("E43").Value2 = 78400
$ws.Range("D47").Value2 = 951800
$ws.Range("E47").Value2 = 937800
$ws.Range("D48").Value2 = 17100
$ws.Range("E48").Value2 = 16600
$ws.Range("D49").Value2 = 104400
$ws.Range("E49").Value2 = 105800
$ws.Range("D52").Value2 = 71400
$ws.Range("E52").Value2 = 65100
$ws.Range("D54").Value2 = 2321400
$ws.Range("E54").Value2 = 2322200
$ws.Range("D57").Value2 = 71000
$ws.Range("E57").Value2 = 52100
$ws.Range("D59").Value2 = 859100
$ws.Range("E59").Value2 = 991500
$ws.Range("D61").Value2 = 160100
$ws.Range("E61").Value2 = 160700
$ws.Range("D66").Value2 = 1801200
$ws.Range("E66").Value2 = 1793400
$ws.Range("D72").Value2 = 140500
$ws.Range("E72").Value2 = 154200
$ws.Range("D76").Value2 = 520200
$ws.Range("E76").Value2 = 528900
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = -11100
$ws.Range("E81").Value2 = -11700
$ws.Range("D83").Value2 = 2300
$ws.Range("E83").Value2 = 3300
$ws.Range("D89").Value2 = -102400
$ws.Range("E89").Value2 = 34000
$ws.Range("D91").Value2 = -1400
$ws.Range("E91").Value2 = -700
$ws.Range("D94").Value2 = -22500
$ws.Range("E94").Value2 = 3600
$ws.Range("D100").Value2 = 3000
$ws.Range("E100").Value2 = 26300
$ws.Range("D102").Value2 = -122000
$ws.Range("E102").Value2 = 63900
